{"js": "// Replace every occurrence of \"John Doe\" with \"Jane Doe\" throughout the\n// document body (the \"Employee Name:\" line and the \"Generated for:\" line).\nconst body = context.document.body;\nconst results = body.search(\"John Doe\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Jane Doe\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace every occurrence of \"John Doe\" with \"Jane Doe\" throughout the\n# document body (the \"Employee Name:\" line and the \"Generated for:\" line).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"John Doe\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Jane Doe\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
